# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the regenerated scrape data (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        7  = 1781
        10 = 152
        11 = 2097
        12 = 28
        13 = 146
        14 = 1347
        15 = 471
        16 = 24
        21 = 38
        24 = 14
        25 = 1133
        27 = 338
        28 = 177
        30 = 333
    }
    "全部类型" = @{
        7  = 1781
        11 = 152
        12 = 2097
        13 = 28
        14 = 146
        15 = 1347
        16 = 471
        17 = 24
        22 = 38
        25 = 14
        26 = 1133
        28 = 338
        29 = 177
        31 = 333
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
